$d = $word.ActiveDocument

# Update the document date/time stamp (day and time) in the "Date" style paragraph.
$d.Content.Find.Execute("19 Feb 2021 13:26:24", $true, $false, $false, $false, $false,
                         $true, 1, $false, "26 Feb 2021 14:49:59", 2)
